$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "sda"
$ws.Range("C3").Value = "ada"
$ws.Range("D3").Value = "ff"

$ws.Range("D3").Select()
